$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remember C2's current (hyperlink) style so we can restore it after the
# collection-wide hyperlink delete below (this runtime's Hyperlinks.Delete()
# clears every hyperlink on the sheet, not just the one it was called on).
$c2Style = $ws.Range("C2").Style
$c2Target = "https://www.upwork.com/nx/search/jobs/?client_hires=1-9,10-&nbs=1&payment_verified=1&q=video%20editing&sort=recency"

# Drop every hyperlink on the sheet, then re-create only the one that
# should survive (C2). The "video editor" hyperlink that lived on C3 is
# gone for good.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C2"), $c2Target)
$ws.Range("C2").Style = $c2Style

# Row 2 (id=1, "video editing"): refreshed RSS feed link (new token/user ids,
# no location filter) - search_url/description stay as they were.
$ws.Range("B2").Value = "https://www.upwork.com/ab/feed/jobs/rss?client_hires=1-9%2C10-&paging=NaN-undefined&payment_verified=1&q=video%20editing&sort=recency&api_params=1&securityToken=418adc6b3d5cfe830ce53f53b359cf3ed0874223d1a0521e8482731da9873c591396925763d482d9d7f4595e2a28196fcf5f1af259c8298d64066d1b522d2fdf&userUid=1795161265690873856&orgUid=1795161265690873857"

# Row 3 (id=2): repurposed from the old "video editor" row into the
# "video editing, us-only" row.
$ws.Range("B3").Value = "https://www.upwork.com/ab/feed/jobs/rss?client_hires=1-9%2C10-&paging=NaN-undefined&payment_verified=1&q=video%20editing&sort=recency&user_location_match=1&api_params=1&securityToken=418adc6b3d5cfe830ce53f53b359cf3ed0874223d1a0521e8482731da9873c591396925763d482d9d7f4595e2a28196fcf5f1af259c8298d64066d1b522d2fdf&userUid=1795161265690873856&orgUid=1795161265690873857"
$ws.Range("C3").Value = "https://www.upwork.com/nx/search/jobs/?client_hires=1-9,10-&nbs=1&payment_verified=1&q=video%20editing&sort=recency&user_location_match=1"
$ws.Range("D3").Value = "payment verified, 1 to 9 and 10+ hires, us-only ""video editor"""

$ws.Range("C3").Select()
